# Modified state names and descriptions per input from 3ric
#
# Renames (applied as exact-match whole-cell text substitutions across the
# used range). The order below matches the order the shared strings were
# appended to the workbook's string table by the original edit:
#   "Steering"                                                  -> "Manual"
#   "Halt"                                                      -> "Armed"
#   "Shore station command to transition to Steering state."  -> "Shore station command to transition to Active state."
#   "Shore station command to transition to Halt state."      -> "Shore station command to transition to Armed state."
#   "Arduino indicates transition to Halt state"               -> "Arduino indicates transition to Armed state"
#   "Start"                                                     -> "Power Up"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$renames = @(
    @("Steering", "Manual"),
    @("Halt", "Armed"),
    @("Shore station command to transition to Steering state.", "Shore station command to transition to Active state."),
    @("Shore station command to transition to Halt state.", "Shore station command to transition to Armed state."),
    @("Arduino indicates transition to Halt state", "Arduino indicates transition to Armed state"),
    @("Start", "Power Up")
)

$usedRange = $ws.UsedRange
$rowCount = $usedRange.Rows.Count
$colCount = $usedRange.Columns.Count

foreach ($pair in $renames) {
    $oldVal = $pair[0]
    $newVal = $pair[1]
    for ($r = 1; $r -le $rowCount; $r++) {
        for ($c = 1; $c -le $colCount; $c++) {
            $cell = $usedRange.Cells.Item($r, $c)
            $val = $cell.Value2
            if ($val -eq $oldVal) {
                $cell.Value = $newVal
            }
        }
    }
}

$ws.Range("A4").Select()
